# Daily attendance processing - 2026-01-23 20:01:23
# Swap the order of the "Recorded By" contributors in column G:
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
# Only touches cells whose text is exactly the old value, leaving any
# single-contributor cells ("dnasr281@gmail.com" or "System" alone) unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Text
    if ($v -eq $oldValue) {
        $cell.Value = $newValue
    }
}
